$d = $word.ActiveDocument

# The text we need to split lives inside a legacy VML textbox (w:pict/v:textbox),
# which is not reachable through the normal Shapes/Paragraphs/Find object model in
# this document (that content isn't part of the main story range). We fall back to
# reading/rewriting the document's OOXML package directly through WordOpenXML,
# which IS exposed on the Document/Range object.

$xml = $d.WordOpenXML

$dash = [char]0x2013

$oldRun = '<w:r><w:t xml:space="preserve"> +++INS $w.startDate +++  ' + $dash + '  +++INS $e.endDate+++</w:t></w:r>'

$newRuns = '<w:r><w:t xml:space="preserve"> +++INS $w.startDate +++  ' + $dash + '  +++INS $</w:t></w:r><w:r><w:t>w</w:t></w:r><w:r><w:t>.endDate+++</w:t></w:r>'

if (-not $xml.Contains($oldRun)) {
    throw "expected run text not found in document"
}

$xml = $xml.Replace($oldRun, $newRuns)

$d.WordOpenXML = $xml
